$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add two new columns P1, Q1 continuing the 0..n sequence ---
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
# Copy the header formatting (bold font, border, centered alignment) from O1 onto the new cells
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Data rows 2-25: flip the I/K and M/O columns, then add new P/Q columns of 2's ---
$ws.Range("I2:I25").Value = 2
$ws.Range("K2:K25").Value = 1
$ws.Range("M2:M25").Value = 2
$ws.Range("O2:O25").Value = 1
$ws.Range("P2:Q25").Value = 2
